$wb = $excel.ActiveWorkbook

# Update shared-string labels used across all sheets (NewStreet/MixedStreet/OldStreet/linear/radial)
$labelMap = @{
    "NewStreet"   = "New street"
    "MixedStreet" = "Mixed street"
    "OldStreet"   = "Old street"
    "linear"      = "Series"
    "radial"      = "Parallel"
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    foreach ($row in 2..6) {
        $cell = $ws.Cells.Item($row, 1)   # column A
        $cur = $cell.Value()
        if ($labelMap.ContainsKey($cur)) {
            $cell.Value = $labelMap[$cur]
        }
    }
}

# --- Sheet 1: "Capacitances kWhperK" ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("B5").Value = 1.506087441489418
$ws.Range("D5").Value = 1021.588666395225
$ws.Range("B6").Value = 1.19112801248157
$ws.Range("D6").Value = 1021.273706966217

# --- Sheet 2: "Downward energy kWh" ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 728.1378986128668
$ws.Range("C2").Value = 728.0536107133894
$ws.Range("D2").Value = 6.625455132925002
$ws.Range("E2").Value = 4.888157839120087
$ws.Range("F2").Value = 732.9417685525095
$ws.Range("B3").Value = 567.3792269942551
$ws.Range("C3").Value = 567.3792269928329
$ws.Range("D3").Value = 6.490536311924996
$ws.Range("E3").Value = 4.795462583408751
$ws.Range("F3").Value = 572.1746895762416
$ws.Range("B4").Value = 413.7448922070051
$ws.Range("C4").Value = 413.58855031634
$ws.Range("D4").Value = 6.456398948808332
$ws.Range("E4").Value = 4.74442293016145
$ws.Range("F4").Value = 418.3329732465015
$ws.Range("B5").Value = 2035.489665051741
$ws.Range("C5").Value = 2035.489665051741
$ws.Range("D5").Value = 28.96880013744834
$ws.Range("E5").Value = 25.87363108599965
$ws.Range("F5").Value = 2061.363296137741
$ws.Range("B6").Value = 2035.489665051741
$ws.Range("C6").Value = 2035.489665051741
$ws.Range("D6").Value = 13.78667678452683
$ws.Range("E6").Value = 12.08802500770862
$ws.Range("F6").Value = 2047.57769005945

# --- Sheet 3: "Upward Energy kWh" ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = 978.3246159645652
$ws.Range("C2").Value = 978.1577230865371
$ws.Range("D2").Value = 7.152293661016666
$ws.Range("E2").Value = 7.06487274975
$ws.Range("F2").Value = 985.222595836287
$ws.Range("B3").Value = 795.3352720803679
$ws.Range("C3").Value = 795.335272079975
$ws.Range("D3").Value = 7.104331046333335
$ws.Range("E3").Value = 7.169742839516664
$ws.Range("F3").Value = 802.5050149194916
$ws.Range("B4").Value = 619.5947373738678
$ws.Range("C4").Value = 619.2833570382555
$ws.Range("D4").Value = 7.174307332708334
$ws.Range("E4").Value = 7.128480531408336
$ws.Range("F4").Value = 626.4118375696638
$ws.Range("B5").Value = 2858.956773389188
$ws.Range("C5").Value = 2858.956773389189
$ws.Range("D5").Value = 30.24237991109465
$ws.Range("E5").Value = 30.18670223978383
$ws.Range("F5").Value = 2889.143475628972
$ws.Range("B6").Value = 2858.956773389188
$ws.Range("C6").Value = 2858.956773389188
$ws.Range("D6").Value = 14.35880447359172
$ws.Range("E6").Value = 14.39645471192019
$ws.Range("F6").Value = 2873.353228101108

# --- Sheet 4: "Max upward power kW" ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = 114.148285396
$ws.Range("C2").Value = 114.25449989113
$ws.Range("D2").Value = 22.75799457469999
$ws.Range("E2").Value = 21.8219291373
$ws.Range("F2").Value = 114.25449989113
$ws.Range("B3").Value = 109.60337307274
$ws.Range("C3").Value = 109.60337307287
$ws.Range("D3").Value = 19.21533658469999
$ws.Range("E3").Value = 20.6092917201
$ws.Range("F3").Value = 109.60337307287
$ws.Range("B4").Value = 107.81158514048
$ws.Range("C4").Value = 107.80596518561
$ws.Range("D4").Value = 17.6087516207
$ws.Range("E4").Value = 19.420590118
$ws.Range("F4").Value = 107.80596518561
$ws.Range("B5").Value = 335.341277757649
$ws.Range("C5").Value = 335.341277757649
$ws.Range("D5").Value = 65.79397546363202
$ws.Range("E5").Value = 64.04184993928159
$ws.Range("F5").Value = 335.341277757649
$ws.Range("B6").Value = 335.341277757649
$ws.Range("C6").Value = 335.341277757649
$ws.Range("D6").Value = 54.29547261951701
$ws.Range("E6").Value = 63.19183769733179
$ws.Range("F6").Value = 335.341277757649

# --- Sheet 5: "Stored energy vs capacitance" ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = 3.347704725050633
$ws.Range("C2").Value = 3.347133638452959
$ws.Range("D2").Value = 17.14747368072926
$ws.Range("E2").Value = 16.93788388392589
$ws.Range("F2").Value = 3.366503805686526
$ws.Range("B3").Value = 2.036326170605284
$ws.Range("C3").Value = 2.036326170604278
$ws.Range("D3").Value = 17.03248432040363
$ws.Range("E3").Value = 17.18930771932724
$ws.Range("F3").Value = 2.052491212645102
$ws.Range("B4").Value = 1.267298663200775
$ws.Range("C4").Value = 1.266661776120776
$ws.Range("D4").Value = 14.56734842434912
$ws.Range("E4").Value = 14.47429763202149
$ws.Range("F4").Value = 1.279952803588751
$ws.Range("B5").Value = 2.802671893800524
$ws.Range("C5").Value = 2.802671893800525
$ws.Range("D5").Value = 20.0800956690715
$ws.Range("E5").Value = 20.04312725025529
$ws.Range("F5").Value = 2.828088809778594
$ws.Range("B6").Value = 2.802671893800524
$ws.Range("C6").Value = 2.802671893800524
$ws.Range("D6").Value = 12.05479538985646
$ws.Range("E6").Value = 12.08640428321968
$ws.Range("F6").Value = 2.813499660768371

# --- Sheet 6: "Response time" ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("B2").Value = 0.8125
$ws.Range("C2").Value = 0.8125
$ws.Range("D2").Value = 0.02777777777777778
$ws.Range("E2").Value = 0.02430555555555556
$ws.Range("F2").Value = 0.8125
$ws.Range("B3").Value = 0.8020833333333334
$ws.Range("C3").Value = 0.8020833333333334
$ws.Range("D3").Value = 0.03125
$ws.Range("E3").Value = 0.02777777777777778
$ws.Range("F3").Value = 0.8020833333333334
$ws.Range("B4").Value = 0.6979166666666667
$ws.Range("C4").Value = 0.6979166666666667
$ws.Range("D4").Value = 0.03472222222222222
$ws.Range("E4").Value = 0.02777777777777778
$ws.Range("F4").Value = 0.6979166666666667
$ws.Range("B5").Value = 1.270833333333333
$ws.Range("C5").Value = 1.270833333333333
$ws.Range("D5").Value = 0.0451388888888889
$ws.Range("E5").Value = 0.03819444444444445
$ws.Range("F5").Value = 1.270833333333333
$ws.Range("B6").Value = 1.270833333333333
$ws.Range("C6").Value = 1.270833333333333
$ws.Range("D6").Value = 0.02430555555555556
$ws.Range("E6").Value = 0.02083333333333333
$ws.Range("F6").Value = 1.270833333333333

# --- Sheet 7: "Energy difference kWh" ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("B2").Value = 250.1867173516957
$ws.Range("C2").Value = 250.104112373142
$ws.Range("D2").Value = 0.5268385280915027
$ws.Range("E2").Value = 2.176714910625378
$ws.Range("F2").Value = 252.2808272837674
$ws.Range("B3").Value = 227.9560450861154
$ws.Range("C3").Value = 227.9560450871386
$ws.Range("D3").Value = 0.6137947344095664
$ws.Range("E3").Value = 2.374280256109159
$ws.Range("F3").Value = 230.3303253432478
$ws.Range("B4").Value = 205.8498451668565
$ws.Range("C4").Value = 205.6948067219164
$ws.Range("D4").Value = 0.7179083839009763
$ws.Range("E4").Value = 2.38405760124715
$ws.Range("F4").Value = 208.0788643231635
$ws.Range("B5").Value = 823.4671083374342
$ws.Range("C5").Value = 823.4671083374633
$ws.Range("D5").Value = 1.273579773649544
$ws.Range("E5").Value = 4.313071153785131
$ws.Range("F5").Value = 827.7801794912484
$ws.Range("B6").Value = 823.4671083374342
$ws.Range("C6").Value = 823.4671083374888
$ws.Range("D6").Value = 0.5721276890690206
$ws.Range("E6").Value = 2.308429704196897
$ws.Range("F6").Value = 825.7755380416856

# --- Sheet 8: "efficciency percent" ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("B2").Value = 74.42702419321479
$ws.Range("C2").Value = 74.43110589724748
$ws.Range("D2").Value = 92.63399191941015
$ws.Range("E2").Value = 69.18960909240955
$ws.Range("F2").Value = 74.39351996688852
$ws.Range("B3").Value = 71.33837098780256
$ws.Range("C3").Value = 71.33837098765898
$ws.Range("D3").Value = 91.36027401709086
$ws.Range("E3").Value = 66.88472222893816
$ws.Range("F3").Value = 71.29858118493253
$ws.Range("B4").Value = 66.77669567691645
$ws.Range("C4").Value = 66.78502588761208
$ws.Range("D4").Value = 89.99334219309334
$ws.Range("E4").Value = 66.55587973610642
$ws.Range("F4").Value = 66.78241823603955
$ws.Range("B5").Value = 71.19693742826449
$ws.Range("C5").Value = 71.19693742826448
$ws.Range("D5").Value = 95.78875809923618
$ws.Range("E5").Value = 85.71201610437265
$ws.Range("F5").Value = 71.34859564869012
$ws.Range("B6").Value = 71.19693742826449
$ws.Range("C6").Value = 71.19693742826449
$ws.Range("D6").Value = 96.01549215969001
$ws.Range("E6").Value = 83.96529040486715
$ws.Range("F6").Value = 71.26091112062491

# --- Sheet 9: "Max downward power kW" ---
$ws = $wb.Worksheets.Item(9)
$ws.Range("B2").Value = 58.59806556420001
$ws.Range("C2").Value = 58.5980655644
$ws.Range("D2").Value = 16.5200310559
$ws.Range("E2").Value = 1.516094787630001
$ws.Range("F2").Value = 58.59806556839999
$ws.Range("B3").Value = 45.1898817245
$ws.Range("C3").Value = 45.18988172469999
$ws.Range("D3").Value = 12.8415761722
$ws.Range("E3").Value = 1.457375900210002
$ws.Range("F3").Value = 45.1898825652
$ws.Range("B4").Value = 31.3573040324
$ws.Range("C4").Value = 31.35730403243
$ws.Range("D4").Value = 9.736506157199997
$ws.Range("E4").Value = 1.455261540899999
$ws.Range("F4").Value = 31.3861885265351
$ws.Range("B5").Value = 212.0886596067572
$ws.Range("C5").Value = 212.0886596067572
$ws.Range("D5").Value = 36.64749990315177
$ws.Range("E5").Value = 7.263342628703831
$ws.Range("F5").Value = 213.1577006153055
$ws.Range("B6").Value = 212.0886596067572
$ws.Range("C6").Value = 212.0886596067571
$ws.Range("D6").Value = 36.25249688483122
$ws.Range("E6").Value = 4.574089969523367
$ws.Range("F6").Value = 212.0891038629994

# --- Sheet 10: "Cost difference euro" ---
$ws = $wb.Worksheets.Item(10)
$ws.Range("B2").Value = 477.9511812611709
$ws.Range("C2").Value = 477.9494983402474
$ws.Range("D2").Value = 6.098616604813287
$ws.Range("E2").Value = 2.7114429284793
$ws.Range("F2").Value = 480.6609412687267
$ws.Range("B3").Value = 339.4231819081451
$ws.Range("C3").Value = 339.4231819057059
$ws.Range("D3").Value = 5.876741577525536
$ws.Range("E3").Value = 2.421182327281713
$ws.Range("F3").Value = 341.8443642329876
$ws.Range("B4").Value = 207.8950470401323
$ws.Range("C4").Value = 207.8937435944281
$ws.Range("D4").Value = 5.738490564899621
$ws.Range("E4").Value = 2.360365328902844
$ws.Range("F4").Value = 210.2541089233309
$ws.Range("B5").Value = 1212.022556714328
$ws.Range("C5").Value = 1212.022556714306
$ws.Range("D5").Value = 27.69522036387207
$ws.Range("E5").Value = 21.56055993229165
$ws.Range("F5").Value = 1233.583116646598
$ws.Range("B6").Value = 1212.022556714328
$ws.Range("C6").Value = 1212.02255671423
$ws.Range("D6").Value = 13.21454909548265
$ws.Range("E6").Value = 9.77959530350563
$ws.Range("F6").Value = 1221.802152017735
